$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new sheet right after Login_Info
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Friend_Request_Management"

# --- Login_Info sheet edits ---
# Remove the sample login row 2 content (keeps A2's Hyperlink style, drops the hyperlink + value)
$ws1.Range("A2").Hyperlinks.Delete()
$ws1.Range("A2").ClearContents()
$ws1.Range("B2").ClearContents()
$ws1.Columns.Item(1).ColumnWidth = 23.0

# --- Friend_Request_Management sheet content ---
# (kept in this order so the shared-string table comes out in the same order as the target)
$ws2.Range("A1").Value = "Cancel Friend Requests"
$ws2.Range("B1").Value = "Accept Friend Requests"
$ws2.Range("A2").Value = 9
$ws2.Range("D1").Value = "Description"
$ws2.Range("D2").Value = "Always write in the column no 2 of the sheet."
$ws2.Range("C1").Value = "Skip FRM"
$ws2.Range("D3").WrapText = $true

$ws2.Columns.Item(1).ColumnWidth = 21.1
$ws2.Columns.Item(2).ColumnWidth = 21.26
$ws2.Columns.Item(3).ColumnWidth = 18.92
$ws2.Columns.Item(4).ColumnWidth = 70.42

# --- Selections / active sheet ---
$ws2.Range("C6").Select()
$ws1.Activate()
$ws1.Range("B7").Select()
